$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 682.9231
$ws.Range("I19").Value = 580
$ws.Range("K19").Value = 580
$ws.Range("M19").Value = -405

$ws.Range("H33").Value = 230
$ws.Range("I33").Value = 232.72728
$ws.Range("K33").Value = 232.72728
$ws.Range("M33").Value = -3.727280000000007

$ws.Range("H111").Value = 1656.8572
$ws.Range("I111").Value = 1551.3334
$ws.Range("J111").Value = 1846.8
$ws.Range("K111").Value = 4654.0002
$ws.Range("L111").Value = 5540.4
$ws.Range("M111").Value = -1587.0002
$ws.Range("N111").Value = -11674.4

$ws.Range("H129").Value = 1209.6
$ws.Range("J129").Value = 1753.6666
$ws.Range("L129").Value = 5260.9998
$ws.Range("N129").Value = -15260.9998

$ws.Range("H137").Value = 26317914
$ws.Range("I137").Value = 47620228
$ws.Range("J137").Value = 3288.4119
$ws.Range("K137").Value = 142860684
$ws.Range("L137").Value = 9865.235700000001
$ws.Range("M137").Value = -142858134
$ws.Range("N137").Value = -14965.2357

$ws.Range("H141").Value = 3577
$ws.Range("I141").Value = 2726.5386
$ws.Range("J141").Value = 9105
$ws.Range("K141").Value = 8179.6158
$ws.Range("L141").Value = 27315
$ws.Range("M141").Value = -2999.6158
$ws.Range("N141").Value = -37675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 72842.92999999999
$ws.Range("I2").Value = 92447.09
$ws.Range("K2").Value = 92447.09
$ws.Range("M2").Value = -92334.09

$ws.Range("H32").Value = 2632.63
$ws.Range("I32").Value = 2663.0918
$ws.Range("J32").Value = 1140
$ws.Range("K32").Value = 2663.0918
$ws.Range("L32").Value = 1140
$ws.Range("M32").Value = -2376.0918
$ws.Range("N32").Value = -1714

$ws.Range("H45").Value = 1087.1923
$ws.Range("I45").Value = 1073.625
$ws.Range("K45").Value = 1073.625
$ws.Range("M45").Value = -696.625

$ws.Range("H97").Value = 11852.223
$ws.Range("I97").Value = 14952.857
$ws.Range("K97").Value = 14952.857
$ws.Range("M97").Value = -14456.857

$ws.Range("H116").Value = 72842.92999999999
$ws.Range("I116").Value = 92447.09
$ws.Range("K116").Value = 92447.09
$ws.Range("M116").Value = -90153.09

$ws.Range("H132").Value = 2668.5334
$ws.Range("I132").Value = 2352.0833
$ws.Range("J132").Value = 3934.3333
$ws.Range("K132").Value = 7056.249899999999
$ws.Range("L132").Value = 11802.9999
$ws.Range("M132").Value = -4526.249899999999
$ws.Range("N132").Value = -16862.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 72842.92999999999
$ws.Range("I3").Value = 92447.09
$ws.Range("K3").Value = 92447.09
$ws.Range("M3").Value = -92333.09

$ws.Range("H94").Value = 2486.5217
$ws.Range("I94").Value = 1895.8823
$ws.Range("J94").Value = 4160
$ws.Range("K94").Value = 1895.8823
$ws.Range("L94").Value = 4160
$ws.Range("M94").Value = -1444.8823
$ws.Range("N94").Value = -5062

$ws.Range("H134").Value = 3200.1052
$ws.Range("I134").Value = 2256.16
$ws.Range("J134").Value = 5015.385
$ws.Range("K134").Value = 6768.48
$ws.Range("L134").Value = 15046.155
$ws.Range("M134").Value = -4233.48
$ws.Range("N134").Value = -20116.155

$ws.Range("H137").Value = 59793.332
$ws.Range("J137").Value = 59793.332
$ws.Range("L137").Value = 59793.332
$ws.Range("N137").Value = -69993.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1306.0358
$ws.Range("I16").Value = 1103.7059
$ws.Range("J16").Value = 1618.7273
$ws.Range("K16").Value = 1103.7059
$ws.Range("L16").Value = 1618.7273
$ws.Range("M16").Value = -816.7058999999999
$ws.Range("N16").Value = -2192.7273

$ws.Range("H31").Value = 5582.525
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5582.525
$ws.Range("K31").Value = 0
$ws.Range("N31").Value = -6172.525
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 5582.525
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5582.525
$ws.Range("K34").Value = 0
$ws.Range("N34").Value = -5986.525
$ws.Range("M34").ClearContents()

$ws.Range("H86").Value = 33335248
$ws.Range("I86").Value = 50001276
$ws.Range("J86").Value = 3195.8
$ws.Range("K86").Value = 50001276
$ws.Range("L86").Value = 3195.8
$ws.Range("M86").Value = -50000153
$ws.Range("N86").Value = -5441.8

$ws.Range("H89").Value = 33335248
$ws.Range("I89").Value = 50001276
$ws.Range("J89").Value = 3195.8
$ws.Range("K89").Value = 250006380
$ws.Range("L89").Value = 15979
$ws.Range("M89").Value = -250000764
$ws.Range("N89").Value = -27211

$ws.Range("H113").Value = 1306.0358
$ws.Range("I113").Value = 1103.7059
$ws.Range("J113").Value = 1618.7273
$ws.Range("K113").Value = 1103.7059
$ws.Range("L113").Value = 1618.7273
$ws.Range("M113").Value = 1066.2941
$ws.Range("N113").Value = -5958.7273

$ws.Range("H127").Value = 34488
$ws.Range("J127").Value = 34488
$ws.Range("L127").Value = 34488
$ws.Range("N127").Value = -44408

$ws.Range("H132").Value = 4168903.5
$ws.Range("I132").Value = 5557271
$ws.Range("J132").Value = 3800.2
$ws.Range("K132").Value = 16671813
$ws.Range("L132").Value = 11400.6
$ws.Range("M132").Value = -16669283
$ws.Range("N132").Value = -16460.6

$ws.Range("H134").Value = 25864120
$ws.Range("I134").Value = 32259508
$ws.Range("J134").Value = 18521266
$ws.Range("K134").Value = 96778524
$ws.Range("L134").Value = 55563798
$ws.Range("M134").Value = -96775989
$ws.Range("N134").Value = -55568868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 7264252.5
$ws.Range("I81").Value = 3313
$ws.Range("K81").Value = 9939
$ws.Range("M81").Value = -8816

$ws.Range("H84").Value = 7264252.5
$ws.Range("I84").Value = 3313
$ws.Range("K84").Value = 29817
$ws.Range("M84").Value = -24201

$ws.Range("H129").Value = 2154
$ws.Range("J129").Value = 1265.1428
$ws.Range("L129").Value = 3795.4284
$ws.Range("N129").Value = -13795.4284

$ws.Range("H131").Value = 8548830
$ws.Range("J131").Value = 9261154
$ws.Range("L131").Value = 27783462
$ws.Range("N131").Value = -27793542

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 0
$ws.Range("N5").Value = -1224
$ws.Range("M5").ClearContents()

$ws.Range("H6").Value = 22636.334
$ws.Range("I6").Value = 3000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2887

$ws.Range("H16").Value = 22636.334
$ws.Range("I16").Value = 3000
$ws.Range("K16").Value = 3000
$ws.Range("M16").Value = -2750

$ws.Range("H80").Value = 40002856
$ws.Range("I80").Value = 2656.0625
$ws.Range("J80").Value = 111114320
$ws.Range("K80").Value = 2656.0625
$ws.Range("L80").Value = 111114320
$ws.Range("M80").Value = -1658.0625
$ws.Range("N80").Value = -111116316

$ws.Range("H83").Value = 40002856
$ws.Range("I83").Value = 2656.0625
$ws.Range("J83").Value = 111114320
$ws.Range("K83").Value = 13280.3125
$ws.Range("L83").Value = 555571600
$ws.Range("M83").Value = -8288.3125
$ws.Range("N83").Value = -555581584

$ws.Range("H122").Value = 2137.8064
$ws.Range("I122").Value = 2104
$ws.Range("J122").Value = 2173.8667
$ws.Range("K122").Value = 6312
$ws.Range("L122").Value = 6521.6001
$ws.Range("M122").Value = -3862
$ws.Range("N122").Value = -11421.6001

$ws.Range("H123").Value = 11851.5
$ws.Range("J123").Value = 11851.5
$ws.Range("L123").Value = 11851.5
$ws.Range("N123").Value = -16751.5

$ws.Range("H132").Value = 3583.16
$ws.Range("I132").Value = 3592.2
$ws.Range("J132").Value = 3569.6
$ws.Range("K132").Value = 10776.6
$ws.Range("L132").Value = 10708.8
$ws.Range("M132").Value = -8246.599999999999
$ws.Range("N132").Value = -15768.8

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2957.4243
$ws.Range("I7").Value = 1832.2222
$ws.Range("J7").Value = 3379.375
$ws.Range("K7").Value = 1832.2222
$ws.Range("L7").Value = 3379.375
$ws.Range("M7").Value = -1720.2222
$ws.Range("N7").Value = -3603.375

$ws.Range("H93").Value = 2246.5881
$ws.Range("I93").Value = 1918.7
$ws.Range("J93").Value = 2715
$ws.Range("K93").Value = 1918.7
$ws.Range("L93").Value = 2715
$ws.Range("M93").Value = -670.7
$ws.Range("N93").Value = -5211

$ws.Range("H126").Value = 2957.4243
$ws.Range("I126").Value = 1832.2222
$ws.Range("J126").Value = 3379.375
$ws.Range("K126").Value = 5496.6666
$ws.Range("L126").Value = 10138.125
$ws.Range("M126").Value = -3026.6666
$ws.Range("N126").Value = -15078.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 745094.75
$ws.Range("I81").Value = 4002331
$ws.Range("K81").Value = 8004662
$ws.Range("M81").Value = -8003601

$ws.Range("H84").Value = 745094.75
$ws.Range("I84").Value = 4002331
$ws.Range("K84").Value = 40023310
$ws.Range("M84").Value = -40018006

$ws.Range("H96").Value = 76924696
$ws.Range("I96").Value = 111112520
$ws.Range("J96").Value = 2085.75
$ws.Range("K96").Value = 111112520
$ws.Range("L96").Value = 2085.75
$ws.Range("M96").Value = -111111147
$ws.Range("N96").Value = -4831.75

$ws.Range("H122").Value = 968.9231
$ws.Range("I122").Value = 925.225
$ws.Range("J122").Value = 1114.5834
$ws.Range("K122").Value = 2775.675
$ws.Range("L122").Value = 3343.7502
$ws.Range("M122").Value = -325.6750000000002
$ws.Range("N122").Value = -8243.7502

$ws.Range("H126").Value = 44168.39
$ws.Range("I126").Value = 48010.906
$ws.Range("K126").Value = 144032.718
$ws.Range("M126").Value = -141562.718

$ws.Range("H136").Value = 6925.375
$ws.Range("I136").Value = 2778.8572
$ws.Range("J136").Value = 10150.444
$ws.Range("K136").Value = 8336.571599999999
$ws.Range("L136").Value = 30451.332
$ws.Range("M136").Value = -5786.571599999999
$ws.Range("N136").Value = -35551.33199999999
